# Fix misspelled unit names on the "Units" sheet (Thief family + Assassin),
# matching the commit's trait/rules cleanup pass, and restore the
# selections that were active when the workbook was saved.

$wb = $excel.ActiveWorkbook

$units = $wb.Worksheets.Item("Units")
$units.Range("A8").Value = "Petty Thief"
$units.Range("A16").Value = "Thief Boss"
$units.Range("A22").Value = "Assassin"

$spells = $wb.Worksheets.Item("Spells & Swift Spells")
$spells.Range("D31").Select()

$units.Activate()
$units.Range("A24").Select()
